{"js": "// Remove the thin grey \"\u2500\u2500\u2500\u2500\" separator paragraphs and the small empty\n// spacer paragraphs (w:spacing w:before=\"40\", no content) that sit\n// between a code-example table and the following Heading3, while\n// leaving every other paragraph - including all paragraphs that hold\n// inline images - untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/spaceBefore\");\nawait context.sync();\n\nconst toDelete = [];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text || \"\";\n  const trimmed = text.trim();\n\n  const isSeparatorLine = trimmed.length > 0 && /^[\u2500]+$/.test(trimmed);\n  // Empty spacer paragraphs of interest use \"w:spacing w:before=40\"\n  // twentieths of a point == 2pt. Paragraphs that merely host an\n  // inline picture also report empty text but use a different\n  // spacing (0 or 4pt before), so checking spaceBefore === 2 keeps\n  // image paragraphs safe.\n  const isSmallEmptySpacer = trimmed.length === 0 && para.spaceBefore === 2;\n\n  if (isSeparatorLine || isSmallEmptySpacer) {\n    toDelete.push(para);\n  }\n}\n\n// Delete from the end backwards so earlier indices/objects stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the thin grey \"\u2500\u2500\u2500\u2500\" separator paragraphs and the small empty\n# spacer paragraphs (w:spacing w:before=\"40\", no content) that sit\n# between a code-example table and the following Heading3, while\n# leaving every other paragraph - including all paragraphs that hold\n# inline images - untouched.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$toDelete = @()\n\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $txt = $p.Range.Text\n  $trimmed = $txt.Trim()\n\n  $isSeparatorLine = ($trimmed.Length -gt 0) -and ($trimmed -match '^[\\x{2500}]+$')\n  # Empty spacer paragraphs of interest use \"w:spacing w:before=40\"\n  # twentieths of a point == 2pt. Paragraphs that merely host an\n  # inline picture also report empty text but use a different\n  # spacing (0 or 4pt before), so checking SpaceBefore -eq 2 keeps\n  # image paragraphs safe.\n  $isSmallEmptySpacer = ($trimmed.Length -eq 0) -and ($p.SpaceBefore -eq 2)\n\n  if ($isSeparatorLine -or $isSmallEmptySpacer) {\n    $toDelete += $i\n  }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$sorted = $toDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n  $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
